$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27 (R43) was updated to reflect a new resistor part (232k instead of 220k).

# PartNumber / Manufacturer Part Number 1
$ws.Range("B27").Value = "ERJ2RKD2323X"
$ws.Range("N27").Value = "ERJ2RKD2323X"

# Name
$ws.Range("C27").Value = "232k"

# Description
$ws.Range("E27").Value = "Resistor - 0402 -  232k - 0.5% - 62mW"

# Link
$ws.Range("G27").Value = "https://cz.mouser.com/ProductDetail/Panasonic/ERJ2RKD2323X?qs=YCa%2FAAYMW02KvIWpqTK%252Bgg%3D%3D"

# Supplier Currency 1 / Supplier Subtotal per Board 1 / Supplier Stock 1 / Supplier Unit Price 1
$ws.Range("I27").Value = "EUR"
$ws.Range("J27").Value = 0.12243
$ws.Range("K27").Value = 21020
$ws.Range("L27").Value = 0.12243
